$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8; this shifts old rows 8..133 down to 9..134
$ws.Rows("8:8").Insert()

# The row that is now at 9 holds the full old row-8 record (all columns A..R).
# Copy that whole record into the newly inserted row 8 so that the "template"
# columns (A,B,C,E,F,G,H,I,O,R) are populated identically, then we will
# overwrite the columns that actually carry new data (D,K,L,M,P).
$ws.Range("A9:R9").Copy()
$ws.Range("A8:R8").PasteSpecial(-4104)  # xlPasteAll
$excel.CutCopyMode = $false

# Now set the new values for row 8 per the update
$ws.Range("D8").Value = 44496
$ws.Range("K8").Value = 7000
$ws.Range("L8").Value = 7000
$ws.Range("M8").Value = 7000
$ws.Range("P8").Value = 1167
